$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 85
$ws.Range("B8").Value = 100
$ws.Range("C8").Value = 99.99998474121094
$ws.Range("D8").Value = 96.99999237060547

$ws.Range("A9").Value = 85
$ws.Range("B9").Value = 100
$ws.Range("C9").Value = 99.99998474121094
$ws.Range("D9").Value = 96.99999237060547
